# Insert a new weekly price record at row 437 for "Feria Lagunitas de Puerto
# Montt" / Piña - Caramelo, pushing the existing rows 437..476 down to
# 438..477 (dimension grows from A1:T476 to A1:T477).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("437:437").Insert()

$ws.Cells.Item(437, 1).Value = 4
$ws.Cells.Item(437, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(437, 3).Value = "Los Lagos"
$ws.Cells.Item(437, 4).Value = 45212
$ws.Cells.Item(437, 5).Value = 10
$ws.Cells.Item(437, 6).Value = "Fruta"
$ws.Cells.Item(437, 7).Value = 100108
$ws.Cells.Item(437, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(437, 9).Value = 100108005
$ws.Cells.Item(437, 10).Value = "Piña"
$ws.Cells.Item(437, 11).Value = "Caramelo"
$ws.Cells.Item(437, 12).Value = "Segunda"
$ws.Cells.Item(437, 13).Value = 100
$ws.Cells.Item(437, 14).Value = 25000
$ws.Cells.Item(437, 15).Value = 25000
$ws.Cells.Item(437, 16).Value = 25000
$ws.Cells.Item(437, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(437, 18).Value = "Ecuador"
$ws.Cells.Item(437, 19).Value = 1786
$ws.Cells.Item(437, 20).Value = 14
